$wb = $excel.ActiveWorkbook

# --- Sheet "Info" (sheet1): update Objetivo/Tiempo result row ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 11191109743103.74
$wsInfo.Range("B2").Value = 0.2069997787475586

# --- Sheet "Activados" (sheet2): change Proceso from 4 to 2, update Tiempo values, add rows 5-8 ---
$wsActivados = $wb.Worksheets.Item("Activados")
$wsActivados.Range("A2").Value = 2
$wsActivados.Range("B2").Value = 0

$wsActivados.Range("A3").Value = 2
$wsActivados.Range("B3").Value = 60

$wsActivados.Range("A4").Value = 2
$wsActivados.Range("B4").Value = 120

$wsActivados.Range("A5").Value = 2
$wsActivados.Range("B5").Value = 180

$wsActivados.Range("A6").Value = 2
$wsActivados.Range("B6").Value = 240

$wsActivados.Range("A7").Value = 2
$wsActivados.Range("B7").Value = 300

$wsActivados.Range("A8").Value = 2
$wsActivados.Range("B8").Value = 360

# --- Sheet "Operando" (sheet3): change Proceso column (A) from 4 to 2 for rows 2-366 ---
$wsOperando = $wb.Worksheets.Item("Operando")
$wsOperando.Range("A2:A366").Value = 2

# --- Sheet "Contaminantes" (sheet6): update Z and Concentracion values ---
$wsContaminantes = $wb.Worksheets.Item("Contaminantes")

$wsContaminantes.Range("B2").Value = 10107185508000
$wsContaminantes.Range("C2").Value = 374.8499999999999

$wsContaminantes.Range("B3").Value = 606673799999.9999
$wsContaminantes.Range("C3").Value = 22.5

$wsContaminantes.Range("B4").Value = 206269092000
$wsContaminantes.Range("C4").Value = 7.649999999999999

$wsContaminantes.Range("B5").Value = 379103.7167999999
$wsContaminantes.Range("C5").Value = 0.00001406

$wsContaminantes.Range("B6").Value = 270980964000
$wsContaminantes.Range("C6").Value = 10.05

Write-Host "Edits applied"
